$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 3522
$ws.Range("B2").Value = 4554
$ws.Range("C2").Value = 4487
$ws.Range("D2").Value = 5816
$ws.Range("E2").Value = 6537
$ws.Range("F2").Value = 5864
